$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet tab to reflect new "through" date
$ws.Name = "Through 2021-10-21"

# Row 9 (July) - 2021 columns T,U,V
$ws.Range("T9").Value = 11
$ws.Range("U9").Value = 138
$ws.Range("V9").Value = 0.0738

# Row 12 (October, through 10-20 -> 10-21)
$ws.Range("A12").Value = "October (through 10-21)"
$ws.Range("F12").Value = 30
$ws.Range("G12").Value = 0.0909
$ws.Range("I12").Value = 35
$ws.Range("J12").Value = 0.1667
$ws.Range("L12").Value = 47
$ws.Range("M12").Value = 0.06
$ws.Range("O12").Value = 30
$ws.Range("P12").Value = 0.1176
$ws.Range("R12").Value = 101
$ws.Range("U12").Value = 129

# Row 13 (Total)
$ws.Range("F13").Value = 413
$ws.Range("G13").Value = 0.1061
$ws.Range("I13").Value = 612
$ws.Range("J13").Value = 0.0852
$ws.Range("L13").Value = 534
$ws.Range("M13").Value = 0.107
$ws.Range("O13").Value = 409
$ws.Range("P13").Value = 0.1031
$ws.Range("R13").Value = 949
$ws.Range("S13").Value = 0.0529
$ws.Range("T13").Value = 82
$ws.Range("U13").Value = 1294
$ws.Range("V13").Value = 0.0596
